$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values: column A becomes a simple row index (1..13),
# column B gets new sample values, and three extra rows (11-13) are added.
$colA = @(1,2,3,4,5,6,7,8,9,10,11,12,13)
$colB = @(412,628,305,342,373,437,574,607,653,667,642,429,390)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# Update the chart: style 6 -> 7, and extend the series references to A1:A13 / B1:B13
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.ChartStyle = 7

$series = $chart.SeriesCollection(1)
$series.XValues = "='Sheet1'!`$A`$1:`$A`$13"
$series.Values = "='Sheet1'!`$B`$1:`$B`$13"
